$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new entry for row 39: Hours = 2, Completed = "Load multiply images on the screen. "
$ws.Range("B39").Value = 2
$ws.Range("C39").Value = "Load multiply images on the screen. "

# Update the SUM formula in B48 to include the new row
$ws.Range("B48").Formula = "=SUM(B2:B39)"

# Update the view state (scrolled position & selection)
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$ws.Range("C39").Select()
